$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview" -- widen the status columns (E, F) to fit the new,
# longer "Handed back: in sync with en-US" status text and refresh the
# status cell values themselves.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -- a handback just completed: the target/source file link
# and the handback file name are now populated, and the handback datetime
# is stamped.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "626b4323-00c9-40aa-a9f5-62b898b4db85.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39979fc2eca5466896d188b359b79b2896e4c134/e2e/626b4323-00c9-40aa-a9f5-62b898b4db85.md", "", "", "626b4323-00c9-40aa-a9f5-62b898b4db85.md")
$wsZh.Range("J2").Value = "626b4323-00c9-40aa-a9f5-62b898b4db85.b717b821591561dbfa949a35910d14bc008fe1f9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-06 08:32:25"

$wsZh.Range("I3").Value = "ba22b89f-e772-4e4a-ae87-564aa235f5e5.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39979fc2eca5466896d188b359b79b2896e4c134/e2e/ba22b89f-e772-4e4a-ae87-564aa235f5e5.md", "", "", "ba22b89f-e772-4e4a-ae87-564aa235f5e5.md")
$wsZh.Range("J3").Value = "ba22b89f-e772-4e4a-ae87-564aa235f5e5.1392f63d34a480156e198ec719111036e083f29b.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-06 08:32:25"

# ---------------------------------------------------------------------------
# Sheet "de-de" -- same handback, for the de-de locale.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "626b4323-00c9-40aa-a9f5-62b898b4db85.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39979fc2eca5466896d188b359b79b2896e4c134/e2e/626b4323-00c9-40aa-a9f5-62b898b4db85.md", "", "", "626b4323-00c9-40aa-a9f5-62b898b4db85.md")
$wsDe.Range("J2").Value = "626b4323-00c9-40aa-a9f5-62b898b4db85.b717b821591561dbfa949a35910d14bc008fe1f9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-06 08:32:33"

$wsDe.Range("I3").Value = "ba22b89f-e772-4e4a-ae87-564aa235f5e5.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/39979fc2eca5466896d188b359b79b2896e4c134/e2e/ba22b89f-e772-4e4a-ae87-564aa235f5e5.md", "", "", "ba22b89f-e772-4e4a-ae87-564aa235f5e5.md")
$wsDe.Range("J3").Value = "ba22b89f-e772-4e4a-ae87-564aa235f5e5.1392f63d34a480156e198ec719111036e083f29b.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-06 08:32:33"

Write-Host "Handback report generated."
